$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This edit regenerates the localization-status report for a new handoff
# cycle ("Generate Report for Handoff"):
#   - 43ad64a5-f232-48e2-909d-d15f4e532309.md moves from "Handed back: in
#     sync with en-US" to "Ready for handoff" (a fresh handoff was queued
#     for it), and its row moves below the 4468ce90 file in every sheet.
#   - 4468ce90-7ace-4bbe-927a-4e931ed44153.md stays "Handed back: in sync
#     with en-US" and now sorts first.
#   - The per-language sheets (zh-cn, de-de) get a new "Latest Handoff
#     Datetime" for 43ad64a5's row, and its "Latest Handback" columns
#     revert to the previous handoff-cycle file/date (a new cycle begins).
# ---------------------------------------------------------------------------

# ============================= Sheet "Overview" =============================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "4468ce90-7ace-4bbe-927a-4e931ed44153.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"

$ws1.Range("A3").Value = "43ad64a5-f232-48e2-909d-d15f4e532309.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

# Recreate the hyperlinks in the same order (same rIds / same target URLs),
# only the displayed text changes to follow the new row contents.
$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/e2e/4468ce90-7ace-4bbe-927a-4e931ed44153.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/.localization-config", "", "", ".localization-config") | Out-Null

# ============================= Sheet "zh-cn" =============================
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 2 -> 4468ce90 file
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-10 18:44:51"
$ws2.Range("F2").Value = "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-10 18:47:33"
$ws2.Range("H2").Value = "Include"

# Row 3 -> 43ad64a5 file, now ready for a new handoff
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-10 18:48:17"
$ws2.Range("F3").Value = "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-03-10 18:47:33"
$ws2.Range("H3").Value = "Include"

# Row 4 (.localization-config) content is unchanged
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

# Recreate hyperlinks in the same order / with the same target URLs, just
# with new display text reflecting the new row contents.
$ws2.Range("A1").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e0985225238a92082d9230478192b995dfb2671/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2833fcdb5ee8eb5daa7557814458bf30d779ac92/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f8f2fad4774db951f94892d9cab2ef41eecd0bf0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/e2e/4468ce90-7ace-4bbe-927a-4e931ed44153.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e0985225238a92082d9230478192b995dfb2671/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.zh-cn.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2833fcdb5ee8eb5daa7557814458bf30d779ac92/e2e/4468ce90-7ace-4bbe-927a-4e931ed44153.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f8f2fad4774db951f94892d9cab2ef41eecd0bf0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.zh-cn.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/.localization-config", "", "", ".localization-config") | Out-Null

# ============================= Sheet "de-de" =============================
$ws3 = $wb.Worksheets.Item("de-de")

# Row 2 -> 4468ce90 file
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-10 18:45:51"
$ws3.Range("F2").Value = "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-10 18:47:44"
$ws3.Range("H2").Value = "Include"

# Row 3 -> 43ad64a5 file, now ready for a new handoff
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-10 18:48:21"
$ws3.Range("F3").Value = "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf"
$ws3.Range("G3").Value = "2016-03-10 18:47:44"
$ws3.Range("H3").Value = "Include"

# Row 4 (.localization-config) content is unchanged
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

# Recreate hyperlinks in the same order / with the same target URLs, just
# with new display text reflecting the new row contents.
$ws3.Range("A1").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4cb8ea8caf8c98b17de2b98af3cfe40f7286e890/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5537a0eb5e40c08714f38f2dc9dc16164b0dac8b/e2e/43ad64a5-f232-48e2-909d-d15f4e532309.md", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/345361376405f7815c42a2044413dd35918da6ce/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf", "", "", "4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/e2e/4468ce90-7ace-4bbe-927a-4e931ed44153.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4cb8ea8caf8c98b17de2b98af3cfe40f7286e890/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.de-de.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5537a0eb5e40c08714f38f2dc9dc16164b0dac8b/e2e/4468ce90-7ace-4bbe-927a-4e931ed44153.md", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/345361376405f7815c42a2044413dd35918da6ce/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4468ce90-7ace-4bbe-927a-4e931ed44153.de700fac5a6640ab97e8c3c2423287d3cc985ba5.de-de.xlf", "", "", "43ad64a5-f232-48e2-909d-d15f4e532309.4633d0e21eabe65683c5fff0e626f3952a175e27.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1304098d8b70371fda361610f3b348be849fb738/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Report regenerated for handoff."
